$wb = $excel.ActiveWorkbook

# ----- Sheet: ALC -----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 154.25
$ws.Range("I33").Value = 147.71428
$ws.Range("K33").Value = 147.71428
$ws.Range("M33").Value = 81.28572
$ws.Range("H64").Value = 3299.6667
$ws.Range("I64").Value = 2999
$ws.Range("J64").Value = 3450
$ws.Range("K64").Value = 2999
$ws.Range("L64").Value = 3450
$ws.Range("M64").Value = -2751
$ws.Range("N64").Value = -3946
$ws.Range("H67").Value = 3299.6667
$ws.Range("I67").Value = 2999
$ws.Range("J67").Value = 3450
$ws.Range("K67").Value = 2999
$ws.Range("L67").Value = 3450
$ws.Range("M67").Value = -2141
$ws.Range("N67").Value = -5166
$ws.Range("H86").Value = 0
$ws.Range("J86").Value = 0
$ws.Range("L86").ClearContents()
$ws.Range("N86").Value = 0
$ws.Range("H89").Value = 0
$ws.Range("J89").Value = 0
$ws.Range("L89").ClearContents()
$ws.Range("N89").Value = 0
$ws.Range("H113").Value = 3985.875
$ws.Range("I113").Value = 3985.875
$ws.Range("K113").Value = 3985.875
$ws.Range("M113").Value = -731.875
$ws.Range("H132").Value = 2665.7666
$ws.Range("I132").Value = 1199.9
$ws.Range("K132").Value = 3599.7
$ws.Range("M132").Value = -1069.7
$ws.Range("H141").Value = 6616.6665
$ws.Range("I141").Value = 4000
$ws.Range("K141").Value = 12000
$ws.Range("M141").Value = -6820

# ----- Sheet: ARM -----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2240.8823
$ws.Range("I45").Value = 1287.8889
$ws.Range("K45").Value = 1287.8889
$ws.Range("M45").Value = -910.8888999999999
$ws.Range("H122").Value = 2108.5557
$ws.Range("I122").Value = 2108.5557
$ws.Range("K122").Value = 6325.6671
$ws.Range("M122").Value = -3875.6671
$ws.Range("H132").Value = 1649.75
$ws.Range("I132").Value = 1649.75
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 4949.25
$ws.Range("L132").Value = 0
$ws.Range("M132").ClearContents()
$ws.Range("N132").Value = -2419.25
$ws.Range("H135").Value = 80000
$ws.Range("J135").Value = 80000
$ws.Range("L135").Value = 80000
$ws.Range("N135").Value = -90140

# ----- Sheet: BSM -----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2435.625
$ws.Range("I20").Value = 1797.2
$ws.Range("K20").Value = 1797.2
$ws.Range("M20").Value = -1550.2
$ws.Range("H22").Value = 4758.5713
$ws.Range("I22").Value = 4758.5713
$ws.Range("K22").Value = 4758.5713
$ws.Range("M22").Value = -4585.5713

# ----- Sheet: CRP -----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5517.1353
$ws.Range("I31").Value = 3176.9167
$ws.Range("K31").Value = 3176.9167
$ws.Range("M31").Value = -2881.9167
$ws.Range("H34").Value = 5517.1353
$ws.Range("I34").Value = 3176.9167
$ws.Range("K34").Value = 3176.9167
$ws.Range("M34").Value = -2974.9167
$ws.Range("H58").Value = 6331.6665
$ws.Range("I58").Value = 5000
$ws.Range("J58").Value = 6997.5
$ws.Range("K58").Value = 5000
$ws.Range("L58").Value = 6997.5
$ws.Range("M58").Value = -4797
$ws.Range("N58").Value = -7403.5
$ws.Range("H105").Value = 1026.5
$ws.Range("I105").Value = 1026.5
$ws.Range("K105").Value = 1026.5
$ws.Range("M105").Value = 720.5
$ws.Range("H132").Value = 1694.125
$ws.Range("I132").Value = 1765
$ws.Range("K132").Value = 5295
$ws.Range("M132").Value = -2765
$ws.Range("H134").Value = 2266.6667
$ws.Range("I134").Value = 2266.6667
$ws.Range("K134").Value = 6800.000100000001
$ws.Range("M134").Value = -4265.000100000001
$ws.Range("H136").Value = 6331.6665
$ws.Range("I136").Value = 5000
$ws.Range("J136").Value = 6997.5
$ws.Range("K136").Value = 15000
$ws.Range("L136").Value = 20992.5
$ws.Range("M136").Value = -12450
$ws.Range("N136").Value = -26092.5
$ws.Range("H140").Value = 65780
$ws.Range("J140").Value = 65780
$ws.Range("L140").Value = 65780
$ws.Range("N140").Value = -76140

# ----- Sheet: CUL -----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H86").Value = 917.1429000000001
$ws.Range("I86").Value = 566.6667
$ws.Range("J86").Value = 1180
$ws.Range("K86").Value = 1700.0001
$ws.Range("L86").Value = 3540
$ws.Range("M86").Value = -514.0001
$ws.Range("N86").Value = -5912
$ws.Range("H89").Value = 917.1429000000001
$ws.Range("I89").Value = 566.6667
$ws.Range("J89").Value = 1180
$ws.Range("K89").Value = 5100.0003
$ws.Range("L89").Value = 10620
$ws.Range("M89").Value = 827.9997000000003
$ws.Range("N89").Value = -22476
$ws.Range("H137").Value = 2742
$ws.Range("J137").Value = 3905
$ws.Range("L137").Value = 11715
$ws.Range("N137").Value = -21915

# ----- Sheet: GSM -----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2557.5
$ws.Range("I102").Value = 1755
$ws.Range("K102").Value = 1755
$ws.Range("M102").Value = -133
$ws.Range("H113").Value = 4244
$ws.Range("I113").Value = 3272.077
$ws.Range("K113").Value = 3272.077
$ws.Range("M113").Value = -1102.077
$ws.Range("H126").Value = 2399.4
$ws.Range("I126").Value = 1999.5
$ws.Range("J126").Value = 2666
$ws.Range("K126").Value = 5998.5
$ws.Range("L126").Value = 7998
$ws.Range("M126").Value = -3528.5
$ws.Range("N126").Value = -12938
$ws.Range("H132").Value = 1553.3334
$ws.Range("I132").Value = 1425.5555
$ws.Range("K132").Value = 4276.666499999999
$ws.Range("M132").Value = -1746.666499999999
$ws.Range("H140").Value = 103172.164
$ws.Range("I140").Value = 279697
$ws.Range("J140").Value = 67867.2
$ws.Range("K140").Value = 279697
$ws.Range("L140").Value = 67867.2
$ws.Range("M140").Value = -274517
$ws.Range("N140").Value = -78227.2

# ----- Sheet: LTW -----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1046.5
$ws.Range("I22").Value = 1259.6666
$ws.Range("J22").Value = 833.3333
$ws.Range("K22").Value = 1259.6666
$ws.Range("L22").Value = 833.3333
$ws.Range("M22").Value = -964.6666
$ws.Range("N22").Value = -1423.3333
$ws.Range("H27").Value = 1046.5
$ws.Range("I27").Value = 1259.6666
$ws.Range("J27").Value = 833.3333
$ws.Range("K27").Value = 1259.6666
$ws.Range("L27").Value = 833.3333
$ws.Range("M27").Value = -1152.6666
$ws.Range("N27").Value = -1047.3333
$ws.Range("H31").Value = 422.5
$ws.Range("I31").Value = 363.33334
$ws.Range("J31").Value = 600
$ws.Range("K31").Value = 363.33334
$ws.Range("L31").Value = 600
$ws.Range("M31").Value = -115.33334
$ws.Range("N31").Value = -1096
$ws.Range("H40").Value = 3207
$ws.Range("I40").Value = 3207
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 3207
$ws.Range("L40").Value = 0
$ws.Range("M40").ClearContents()
$ws.Range("N40").Value = -3071
$ws.Range("H122").Value = 2832.1667
$ws.Range("I122").Value = 2832.1667
$ws.Range("K122").Value = 8496.500100000001
$ws.Range("M122").Value = -6046.500100000001

# ----- Sheet: WVR -----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H86").Value = 70325
$ws.Range("J86").Value = 70325
$ws.Range("L86").Value = 70325
$ws.Range("N86").Value = -72571
$ws.Range("H89").Value = 70325
$ws.Range("J89").Value = 70325
$ws.Range("L89").Value = 351625
$ws.Range("N89").Value = -362857
$ws.Range("H122").Value = 3201
$ws.Range("I122").Value = 3201
$ws.Range("K122").Value = 9603
$ws.Range("M122").Value = -7153
